$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:51:31.315157",
    "2021-10-05 10:51:31.315171",
    "2021-10-05 10:51:31.315175",
    "2021-10-05 10:51:31.315178",
    "2021-10-05 10:51:31.315182",
    "2021-10-05 10:51:31.315185",
    "2021-10-05 10:51:31.315188",
    "2021-10-05 10:51:31.315191",
    "2021-10-05 10:51:31.315195",
    "2021-10-05 10:51:31.315198",
    "2021-10-05 10:51:31.315201",
    "2021-10-05 10:51:31.315204",
    "2021-10-05 10:51:31.315207",
    "2021-10-05 10:51:31.315210",
    "2021-10-05 10:51:31.315213",
    "2021-10-05 10:51:31.315217",
    "2021-10-05 10:51:31.315220"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
